$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in the header (F1)
$ws.Range("F1").Value = "Last status check on: 26.02.2022 11:45"

# Row 3 (Tesco) gets refreshed prices/delta/date values
# B3: new current price
$ws.Range("B3").Value = 37.29
# C3: old price (previous B3 value)
$ws.Range("C3").Value = 36.9

# D3: delta, stored as text "+0.39" (force text type without leaving a
# lingering custom number-format style on the cell)
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "+0.39"
$ws.Range("D3").Style = "Normal"

# E3: timestamp stored as plain text instead of a date serial number
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2022-02-26 11:45:08"
$ws.Range("E3").Style = "Normal"
